$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at position 38, shifting rows 38:150 down to 39:151
$ws.Rows("38:38").Insert()

# Populate the newly inserted row's September Details / September Date cells
$ws.Range("R38").Value = "swiggy refunded"
$ws.Range("S38").Value = "2024-09-15 21:05:07"
